$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.267.70"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "3.325.55"
$ws.Range("E3").Value = "  +1.96%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'578.05"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").Value = "'184.66"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").Value = "'6.65"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "3.900.21"
$ws.Range("E12").Value = "  +1.94%  "
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").Value = "'27.31"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "67.466.42"
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("D16").Value = "'0.0000167"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").Value = "3.320.10"
$ws.Range("E17").Value = "  +1.93%  "
$ws.Range("D18").Value = "'443.37"
$ws.Range("E18").Value = "  +6.24%  "
$ws.Range("D19").Value = "'13.54"
$ws.Range("D20").Value = "'5.66"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "'7.69"
$ws.Range("E21").Value = "  +2.29%  "
$ws.Range("D22").Value = "'73.89"
$ws.Range("E22").Value = "  +3.98%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "3.463.76"
$ws.Range("E24").Value = "  +1.93%  "
$ws.Range("D25").Value = "'0.511"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  +1.28%  "
$ws.Range("E27").Value = "  +1.94%  "
$ws.Range("D28").Value = "'9.03"
$ws.Range("E28").Value = "  -3.53%  "
$ws.Range("D29").Value = "'0.985"
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("D30").Value = "'1.97"
$ws.Range("E30").Value = "  +1.07%  "
$ws.Range("D31").Value = "'22.85"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("D32").Value = "'5.32"
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "'6.79"
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("E35").Value = "  -0.88%  "
$ws.Range("D36").Value = "'1.50"
$ws.Range("E36").Value = "  +4.29%  "
$ws.Range("D37").Value = "'161.50"
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("D40").Value = "2.790.42"
$ws.Range("E40").Value = "  +5.98%  "
$ws.Range("D41").Value = "'0.789"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "'6.21"
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("D45").Value = "'0.0670"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").Value = "'24.57"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("D47").Value = "'2.37"
$ws.Range("E47").Value = "  -2.21%  "
$ws.Range("D48").Value = "'325.03"
$ws.Range("E48").Value = "  -3.47%  "
$ws.Range("D49").Value = "'0.0272"
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("D50").Value = "'0.983"
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "'30.96"
$ws.Range("E51").Value = "  +1.58%  "
